$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date field text on the
#    slide master and every slide layout (9/27/2018 -> 12/7/2018).
# ------------------------------------------------------------------
$master = $p.SlideMaster

for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $shp = $master.Shapes.Item($si)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "12/7/2018"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12/7/2018"
        }
    }
}

# ------------------------------------------------------------------
# 2) Swap the "Karma/Jasmine" testing-framework label for "Jest" on
#    slide 1, touching only that run so the sibling "TestNG," run
#    (and all its formatting) is left untouched.
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text.Contains("Karma/Jasmine")) {
            $target = $shp
        }
    }
}

$tr = $target.TextFrame.TextRange
$fullText = $tr.Text
$startPos = $fullText.IndexOf("Karma/Jasmine") + 1
$old = $tr.Characters($startPos, 13)
$old.Text = "Jest"
